$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.703.89"
$ws.Range("E2").Value = "  +1.94%  "
$ws.Range("D3").Value = "1.900.78"
$ws.Range("E3").Value = "  +2.96%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.27"
$ws.Range("E5").Value = "  +1.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4825"
$ws.Range("E7").Value = "  +1.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2850"
$ws.Range("E8").Value = "  +1.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06562"
$ws.Range("E9").Value = "  +1.44%  "
$ws.Range("D10").Value = "1.951.49"
$ws.Range("E10").Value = "  +5.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07460"
$ws.Range("E11").Value = "  +2.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.72"
$ws.Range("E12").Value = "  +2.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.110"
$ws.Range("E13").Value = "  -0.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.18"
$ws.Range("E14").Value = "  +1.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6683"
$ws.Range("E15").Value = "  +3.72%  "
$ws.Range("D16").Value = "30.686.50"
$ws.Range("E16").Value = "  +2.05%  "
$ws.Range("B17").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C17").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D17").Value = "2.263.86"
$ws.Range("E17").Value = "  +8.04%  "
$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.33"
$ws.Range("E18").Value = "  +0.92%  "
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9999"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007619"
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "232.10"
$ws.Range("E21").Value = "  +4.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.293"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.236"
$ws.Range("E24").Value = "  +2.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "170.07"
$ws.Range("E25").Value = "  +3.98%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.339"
$ws.Range("E26").Value = "  +1.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.79"
$ws.Range("E27").Value = "  +1.80%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.973"
$ws.Range("E28").Value = "  +3.07%  "
$ws.Range("E29").Value = "  -1.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.1013"
$ws.Range("E30").Value = "  +10.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.367"
$ws.Range("E31").Value = "  +3.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.039"
$ws.Range("E32").Value = "  +2.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05120"
$ws.Range("E33").Value = "  +2.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.222"
$ws.Range("E34").Value = "  +7.71%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7602"
$ws.Range("E35").Value = "  +2.98%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.703"
$ws.Range("E36").Value = "  +0.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01895"
$ws.Range("E37").Value = "  +4.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.661"
$ws.Range("E38").Value = "  +2.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9223"
$ws.Range("E39").Value = "  +2.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.086"
$ws.Range("E40").Value = "  +1.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "107.15"
$ws.Range("E41").Value = "  +0.41%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4310"
$ws.Range("E42").Value = "  +1.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.746"
$ws.Range("E44").Value = "  -3.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.441"
$ws.Range("E45").Value = "  +0.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.73"
$ws.Range("E46").Value = "  +1.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1277"
$ws.Range("E47").Value = "  -2.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.489"
$ws.Range("E48").Value = "  -4.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.966"
$ws.Range("E49").Value = "  +2.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.92"
$ws.Range("E50").Value = "  -1.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05681"
$ws.Range("E51").Value = "  +0.26%  "
